# ADD results from server
# Update computed result values on the "2025", "2030", and "2035" sheets
# with refreshed figures coming from the server-side model run.

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 43865.8098096851
$ws.Range("B2").Value = 18877.50563722889
$ws.Range("E2").Value = 135823.52998632
$ws.Range("G2").Value = 42315.16049511674
$ws.Range("H2").Value = 493223.6395174918
$ws.Range("I2").Value = 403072.1206296
$ws.Range("N2").Value = 51593.46271940265
$ws.Range("O2").Value = 69302.79012968208

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 51440.97267772973
$ws.Range("B2").Value = 131602.3685301121
$ws.Range("E2").Value = 112231.8351973442
$ws.Range("H2").Value = 105153.7033457349
$ws.Range("I2").Value = 335356.0663069373
$ws.Range("N2").Value = 24263.37816951394
$ws.Range("O2").Value = 52593.17366863995

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 85834.0918338597
$ws.Range("B2").Value = 2.418741552752816
$ws.Range("I2").Value = 176628.7208800473
$ws.Range("M2").Value = 28128.83060279026
$ws.Range("N2").Value = 11927.45854640803
$ws.Range("O2").Value = 57190.46047809131
